$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 995.2963
$ws.Range("J17").Value = 1107.8636
$ws.Range("L17").Value = 3323.5908
$ws.Range("N17").Value = -3659.5908

$ws.Range("H100").Value = 2009.1818
$ws.Range("I100").Value = 1665.8334
$ws.Range("J100").Value = 2421.2
$ws.Range("K100").Value = 1665.8334
$ws.Range("L100").Value = 2421.2
$ws.Range("M100").Value = -1124.8334
$ws.Range("N100").Value = -3503.2

$ws.Range("H111").Value = 1259
$ws.Range("I111").Value = 1259
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3777
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -710
$ws.Range("N111").ClearContents()

$ws.Range("H113").Value = 3530
$ws.Range("I113").Value = 3530
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3530
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -276
$ws.Range("N113").ClearContents()

$ws.Range("H137").Value = 1582.909
$ws.Range("I137").Value = 2239.4
$ws.Range("J137").Value = 1207.7715
$ws.Range("K137").Value = 6718.200000000001
$ws.Range("L137").Value = 3623.3145
$ws.Range("M137").Value = -4168.200000000001
$ws.Range("N137").Value = -8723.3145

$ws.Range("H138").Value = 2773.3735
$ws.Range("I138").Value = 2445.5881
$ws.Range("J138").Value = 2857.803
$ws.Range("K138").Value = 7336.7643
$ws.Range("L138").Value = 8573.409
$ws.Range("M138").Value = -2196.7643
$ws.Range("N138").Value = -18853.409

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14781.3
$ws.Range("I32").Value = 10098.964
$ws.Range("J32").Value = 36111.945
$ws.Range("K32").Value = 10098.964
$ws.Range("L32").Value = 36111.945
$ws.Range("M32").Value = -9811.964
$ws.Range("N32").Value = -36685.945

$ws.Range("H74").Value = 1340.2903
$ws.Range("I74").Value = 1184.7826
$ws.Range("J74").Value = 1787.375
$ws.Range("K74").Value = 1184.7826
$ws.Range("L74").Value = 1787.375
$ws.Range("M74").Value = -310.7826
$ws.Range("N74").Value = -3535.375

$ws.Range("H77").Value = 1340.2903
$ws.Range("I77").Value = 1184.7826
$ws.Range("J77").Value = 1787.375
$ws.Range("K77").Value = 5923.913
$ws.Range("L77").Value = 8936.875
$ws.Range("M77").Value = -1555.913
$ws.Range("N77").Value = -17672.875

$ws.Range("H95").Value = 12547.066
$ws.Range("J95").Value = 12547.066
$ws.Range("L95").Value = 12547.066
$ws.Range("N95").Value = -18039.066

$ws.Range("H101").Value = 12692.556
$ws.Range("J101").Value = 12692.556
$ws.Range("L101").Value = 12692.556
$ws.Range("N101").Value = -19182.556

$ws.Range("H122").Value = 1658.2632
$ws.Range("I122").Value = 1755.7273
$ws.Range("J122").Value = 1524.25
$ws.Range("K122").Value = 5267.1819
$ws.Range("L122").Value = 4572.75
$ws.Range("M122").Value = -2817.1819
$ws.Range("N122").Value = -9472.75

$ws.Range("H132").Value = 2723791.8
$ws.Range("I132").Value = 6535.5454
$ws.Range("J132").Value = 5214610
$ws.Range("K132").Value = 19606.6362
$ws.Range("L132").Value = 15643830
$ws.Range("M132").Value = -17076.6362
$ws.Range("N132").Value = -15648890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1592.4286
$ws.Range("I20").Value = 1510.5555
$ws.Range("J20").Value = 1739.8
$ws.Range("K20").Value = 1510.5555
$ws.Range("L20").Value = 1739.8
$ws.Range("M20").Value = -1263.5555
$ws.Range("N20").Value = -2233.8

$ws.Range("H80").Value = 372.1
$ws.Range("I80").Value = 327.72726
$ws.Range("J80").Value = 426.33334
$ws.Range("K80").Value = 327.72726
$ws.Range("L80").Value = 426.33334
$ws.Range("M80").Value = 670.27274
$ws.Range("N80").Value = -2422.33334

$ws.Range("H83").Value = 372.1
$ws.Range("I83").Value = 327.72726
$ws.Range("J83").Value = 426.33334
$ws.Range("K83").Value = 1638.6363
$ws.Range("L83").Value = 2131.6667
$ws.Range("M83").Value = 3353.3637
$ws.Range("N83").Value = -12115.6667

$ws.Range("H134").Value = 2390.6667
$ws.Range("I134").Value = 1573.3939
$ws.Range("J134").Value = 3889
$ws.Range("K134").Value = 4720.1817
$ws.Range("L134").Value = 11667
$ws.Range("M134").Value = -2185.1817
$ws.Range("N134").Value = -16737

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9982.5
$ws.Range("J4").Value = 9982.5
$ws.Range("L4").Value = 9982.5
$ws.Range("N4").Value = -10206.5

$ws.Range("H43").Value = 21161.889
$ws.Range("J43").Value = 21161.889
$ws.Range("L43").Value = 21161.889
$ws.Range("N43").Value = -21529.889

$ws.Range("H101").Value = 21161.889
$ws.Range("J101").Value = 21161.889
$ws.Range("L101").Value = 21161.889
$ws.Range("N101").Value = -27651.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 610060
$ws.Range("I4").Value = 610060
$ws.Range("K4").Value = 1830180
$ws.Range("M4").Value = -1830068

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 566.3333
$ws.Range("I5").Value = 350
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 350
$ws.Range("L5").Value = 999
$ws.Range("M5").Value = -238
$ws.Range("N5").Value = -1223

$ws.Range("H31").Value = 1543.6666
$ws.Range("I31").Value = 1543.6666
$ws.Range("K31").Value = 1543.6666
$ws.Range("M31").Value = -1251.6666

$ws.Range("H37").Value = 1543.6666
$ws.Range("I37").Value = 1543.6666
$ws.Range("K37").Value = 1543.6666
$ws.Range("M37").Value = -1266.6666

$ws.Range("H70").Value = 8090.5884
$ws.Range("I70").Value = 9971.429
$ws.Range("J70").Value = 6774
$ws.Range("K70").Value = 9971.429
$ws.Range("L70").Value = 6774
$ws.Range("M70").Value = -9701.429
$ws.Range("N70").Value = -7314

$ws.Range("H73").Value = 8090.5884
$ws.Range("I73").Value = 9971.429
$ws.Range("J73").Value = 6774
$ws.Range("K73").Value = 9971.429
$ws.Range("L73").Value = 6774
$ws.Range("M73").Value = -9035.429
$ws.Range("N73").Value = -8646

$ws.Range("H92").Value = 5333.8184
$ws.Range("J92").Value = 5333.8184
$ws.Range("L92").Value = 5333.8184
$ws.Range("N92").Value = -9077.8184

$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws.Range("H126").Value = 4100427.8
$ws.Range("I126").Value = 7814243.5
$ws.Range("K126").Value = 23442730.5
$ws.Range("M126").Value = -23440260.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2503955.8
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 2731542.5
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 2731542.5
$ws.Range("M2").Value = -388
$ws.Range("N2").Value = -2731766.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 151485710
$ws.Range("I2").Value = 30000000
$ws.Range("K2").Value = 30000000
$ws.Range("M2").Value = -29999888

$ws.Range("H68").Value = 27333
$ws.Range("J68").Value = 27333
$ws.Range("L68").Value = 27333
$ws.Range("N68").Value = -28955

$ws.Range("H71").Value = 27333
$ws.Range("J71").Value = 27333
$ws.Range("L71").Value = 81999
$ws.Range("N71").Value = -90111

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H107").Value = 639.86664
$ws.Range("I107").Value = 463.45456
$ws.Range("J107").Value = 1125
$ws.Range("K107").Value = 1390.36368
$ws.Range("L107").Value = 3375
$ws.Range("M107").Value = 529.6363200000001
$ws.Range("N107").Value = -7215

$ws.Range("H117").Value = 33196.2
$ws.Range("J117").Value = 33196.2
$ws.Range("L117").Value = 33196.2
$ws.Range("N117").Value = -42374.2
